# Port battle calculator correction
# - Reorders/corrects the ship names and BR values on the
#   "Shallow water port" sheet (rows 4-21).
# - Fixes the BR-total SUM() ranges on both sheets so they cover
#   the full data range (rows already existed but the totals
#   formula was short by one row).

$wb = $excel.ActiveWorkbook

$deep = $wb.Worksheets.Item("Deep water port")
$shallow = $wb.Worksheets.Item("Shallow water port")

# --- Deep water port: extend BR-total SUM formulas to row 35 ---
$deep.Range("D3").Formula = "=SUM(D4:D35)"
$deep.Range("E3").Formula = "=SUM(E4:E35)"

# --- Shallow water port: extend BR-total SUM formulas to row 21 ---
$shallow.Range("D3").Formula = "=SUM(D4:D21)"
$shallow.Range("E3").Formula = "=SUM(E4:E21)"

# --- Shallow water port: corrected ship / BR listing for rows 4-21 ---
$ships = @(
    @{ Row = 4;  Name = "Hercules";              BR = 100 },
    @{ Row = 5;  Name = "Pandora";                BR = 100 },
    @{ Row = 6;  Name = "Mercury";                BR = 80 },
    @{ Row = 7;  Name = "Mortar Brig";            BR = 80 },
    @{ Row = 8;  Name = "NavyBrig";               BR = 80 },
    @{ Row = 9;  Name = "Niagara";                BR = 80 },
    @{ Row = 10; Name = "Prince de Neufchatel";   BR = 80 },
    @{ Row = 11; Name = "Rattlesnake";            BR = 80 },
    @{ Row = 12; Name = "Rattlesnake Heavy";      BR = 80 },
    @{ Row = 13; Name = "Snow";                   BR = 80 },
    @{ Row = 14; Name = "Brig";                   BR = 70 },
    @{ Row = 15; Name = "Pickle";                 BR = 55 },
    @{ Row = 16; Name = "Cutter";                 BR = 50 },
    @{ Row = 17; Name = "GunBoat";                BR = 50 },
    @{ Row = 18; Name = "Lynx";                   BR = 50 },
    @{ Row = 19; Name = "Privateer";              BR = 50 },
    @{ Row = 20; Name = "Yacht";                  BR = 50 },
    @{ Row = 21; Name = "Yacht Silver";           BR = 50 }
)

foreach ($ship in $ships) {
    $shallow.Cells.Item($ship.Row, 2).Value = $ship.Name
    $shallow.Cells.Item($ship.Row, 3).Value = $ship.BR
}
